$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.154094457626343
$ws.Range("B1").Value = 2.639034509658813
$ws.Range("C1").Value = 2.801393270492554
$ws.Range("D1").Value = 3.401887893676758
$ws.Range("E1").Value = 2.042211771011353
